$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 102: update B102, clear the (empty inline string) C102 cell
$ws.Range("B102").Value = 7.171333983999999
$ws.Range("C102").ClearContents()

# Rows 103-108: A=0, B=7.171333983999999
for ($r = 103; $r -le 108; $r++) {
    $ws.Cells.Item($r, 1).Value = 0
    $ws.Cells.Item($r, 2).Value = 7.171333983999999
}

# Row 109: A=0, B=7.55965918, C109 empty inline string cell
$ws.Cells.Item(109, 1).Value = 0
$ws.Cells.Item(109, 2).Value = 7.55965918
$ws.Cells.Item(109, 3).Value = ""
